# Automatic update of files.
# Bump the "Förändrad" date (column C) from 2025-04-26 (45773) to 2025-04-27 (45774)
# for every data row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45773) {
        $cell.Value2 = 45774
    }
}
